$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "29.135.18"
Set-TextValue $ws.Range("E2") "  -1.02%  "
Set-TextValue $ws.Range("D3") "1.973.71"
Set-TextValue $ws.Range("E3") "  -0.63%  "
Set-TextValue $ws.Range("D4") "1.012"
Set-TextValue $ws.Range("E4") "  +0.41%  "
Set-TextValue $ws.Range("D5") "329.30"
Set-TextValue $ws.Range("E5") "  -0.18%  "
Set-TextValue $ws.Range("D6") "1.011"
Set-TextValue $ws.Range("E6") "  +0.49%  "
Set-TextValue $ws.Range("D7") "0.4957"
Set-TextValue $ws.Range("E7") "  -0.15%  "
Set-TextValue $ws.Range("D8") "0.4210"
Set-TextValue $ws.Range("E8") "  +0.38%  "
Set-TextValue $ws.Range("D9") "54.09"
Set-TextValue $ws.Range("E9") "  +4.22%  "
Set-TextValue $ws.Range("D10") "0.09310"
Set-TextValue $ws.Range("E10") "  +5.22%  "
Set-TextValue $ws.Range("E11") "  -1.79%  "
Set-TextValue $ws.Range("D12") "22.81"
Set-TextValue $ws.Range("E12") "  -2.01%  "
Set-TextValue $ws.Range("D13") "1.981.92"
Set-TextValue $ws.Range("E13") "  +4.40%  "
Set-TextValue $ws.Range("E14") "  -1.71%  "
Set-TextValue $ws.Range("D15") "6.457"
Set-TextValue $ws.Range("E15") "  -0.47%  "
Set-TextValue $ws.Range("D16") "1.014"
Set-TextValue $ws.Range("E16") "  +0.71%  "
Set-TextValue $ws.Range("D17") "0.00001112"
Set-TextValue $ws.Range("E17") "  +0.77%  "
Set-TextValue $ws.Range("D18") "91.84"
Set-TextValue $ws.Range("E18") "  -4.33%  "
Set-TextValue $ws.Range("D19") "0.06724"
Set-TextValue $ws.Range("E19") "  +1.46%  "
Set-TextValue $ws.Range("D20") "19.18"
Set-TextValue $ws.Range("E20") "  -2.54%  "
Set-TextValue $ws.Range("D21") "1.011"
Set-TextValue $ws.Range("E21") "  +0.60%  "
Set-TextValue $ws.Range("D22") "5.956"
Set-TextValue $ws.Range("E22") "  +0.20%  "
Set-TextValue $ws.Range("D23") "29.165.26"
Set-TextValue $ws.Range("E23") "  -0.91%  "
Set-TextValue $ws.Range("E24") "  +1.10%  "
Set-TextValue $ws.Range("D25") "2.264"
Set-TextValue $ws.Range("E25") "  -0.81%  "
Set-TextValue $ws.Range("D26") "2.211.64"
Set-TextValue $ws.Range("E26") "  +2.75%  "
Set-TextValue $ws.Range("D27") "20.78"
Set-TextValue $ws.Range("E27") "  +1.26%  "
Set-TextValue $ws.Range("D28") "156.32"
Set-TextValue $ws.Range("E28") "  -0.78%  "
Set-TextValue $ws.Range("D29") "6.256"
Set-TextValue $ws.Range("E29") "  -3.83%  "
Set-TextValue $ws.Range("E30") "  -2.76%  "
Set-TextValue $ws.Range("E31") "  -0.39%  "
Set-TextValue $ws.Range("E32") "  -0.39%  "
Set-TextValue $ws.Range("D33") "0.09839"
Set-TextValue $ws.Range("E33") "  -0.72%  "
Set-TextValue $ws.Range("D34") "1.503"
Set-TextValue $ws.Range("E34") "  -3.98%  "
Set-TextValue $ws.Range("D35") "5.814"
Set-TextValue $ws.Range("E35") "  -0.46%  "
Set-TextValue $ws.Range("D36") "3.738"
Set-TextValue $ws.Range("E36") "  -1.32%  "
Set-TextValue $ws.Range("D37") "0.02422"
Set-TextValue $ws.Range("E37") "  -0.85%  "
Set-TextValue $ws.Range("D38") "1.328"
Set-TextValue $ws.Range("E38") "  +3.33%  "
Set-TextValue $ws.Range("D39") "0.06425"
Set-TextValue $ws.Range("E39") "  +1.35%  "
Set-TextValue $ws.Range("D40") "9.049"
Set-TextValue $ws.Range("E40") "  -5.12%  "
Set-TextValue $ws.Range("D41") "0.6477"
Set-TextValue $ws.Range("E41") "  -0.29%  "
Set-TextValue $ws.Range("D42") "11.50"
Set-TextValue $ws.Range("E42") "  -1.91%  "
Set-TextValue $ws.Range("D43") "0.2003"
Set-TextValue $ws.Range("E43") "  -2.91%  "
Set-TextValue $ws.Range("D44") "1.011"
Set-TextValue $ws.Range("E44") "  +0.56%  "
Set-TextValue $ws.Range("E45") "  -1.63%  "
Set-TextValue $ws.Range("D46") "1.359"
Set-TextValue $ws.Range("E46") "  +7.82%  "
Set-TextValue $ws.Range("B47") "NEARProtocol"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D47") "2.181"
Set-TextValue $ws.Range("E47") "  -1.27%  "
Set-TextValue $ws.Range("B48") "EnergySwap"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "13.25"
Set-TextValue $ws.Range("E48") "  -1.22%  "
Set-TextValue $ws.Range("E49") "  -1.28%  "
Set-TextValue $ws.Range("E50") "  +2.17%  "
Set-TextValue $ws.Range("D51") "0.06970"
Set-TextValue $ws.Range("E51") "  -0.17%  "
